$d = $word.ActiveDocument

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>4. Content Creation and Media</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Responsibilities:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Produce high-quality video, streaming content, and digital media to engage fans and build the brand’s presence online.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Videographers</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Made a vlog to give an attraction to the supporter to keep support player such as Vlog </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Holiday ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Mic Check , Introduce the player and so on.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Editors</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Give an effort to attract the audience in sound effect, made a good quality, made research about a trend. It became famous some of fans clips the video.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$n = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($n)
$rng = $target.Range
$rng.Collapse(1)
[void]$rng.InsertXML($xml)

Write-Output ("Paragraphs after edit: {0}" -f $d.Paragraphs.Count)
